# Update: Actualización desde MV -datos-
# Appends one new quarterly data row (01-07-2021) to the trade series table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Column A holds a date-like label ("01-07-2021") that must be stored as
# literal text (matching every other "Serie" cell in column A), not parsed
# into a date serial number. Mark the cell as Text before assigning the
# value, then restore the default (unstyled) look used by the rest of the
# data rows, same as cell A2.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item(2, 1).Style

$ws.Cells.Item($row, 2).Value = 23194
$ws.Cells.Item($row, 3).Value = 14700
$ws.Cells.Item($row, 4).Value = 13124
$ws.Cells.Item($row, 5).Value = 7365
$ws.Cells.Item($row, 6).Value = 1129
$ws.Cells.Item($row, 7).Value = 22505
$ws.Cells.Item($row, 8).Value = 7557
$ws.Cells.Item($row, 9).Value = 3026
$ws.Cells.Item($row, 10).Value = 8892
$ws.Cells.Item($row, 11).Value = 3706
$ws.Cells.Item($row, 12).Value = 4607
